$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table into a new column P (year 2022), mirroring the
# formatting already used by the neighbouring column O.

# Row 3 (thin separator row) - new blank, formatted cell at P3.
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

# Row 4 (year header row) - new header cell P4 = 2022.
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2022

# Row 5 (data row) - update the last three existing years and add 2022.
$ws.Range("M5").Value = 2.6
$ws.Range("N5").Value = 2.4
$ws.Range("O5").Value = 3.3

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value = 2.6

$excel.CutCopyMode = 0

# Move the active selection to P3, matching the new selection state.
$ws.Range("P3").Select() | Out-Null
